# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing text storage so that
# numeric-looking strings (e.g. "238.09") are preserved verbatim instead
# of being coerced into floating point numbers by Excel.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '43.696.52'
$ws.Range("E2").Value = '  -0.43%  '

Set-TextValue $ws.Range("D3") '2.334.39'
$ws.Range("E3").Value = '  -1.21%  '

$ws.Range("E4").Value = '  +0.01%  '

Set-TextValue $ws.Range("D5") '238.09'
$ws.Range("E5").Value = '  -1.35%  '

Set-TextValue $ws.Range("D6") '0.661'
$ws.Range("E6").Value = '  -4.04%  '

Set-TextValue $ws.Range("D7") '71.47'
$ws.Range("E7").Value = '  -6.35%  '

$ws.Range("E8").Value = '  -0.07%  '

Set-TextValue $ws.Range("D9") '0.582'
$ws.Range("E9").Value = '  -7.21%  '

Set-TextValue $ws.Range("D10") '0.0982'
$ws.Range("E10").Value = '  -4.15%  '

Set-TextValue $ws.Range("D11") '57.82'
$ws.Range("E11").Value = '  +1.09%  '

Set-TextValue $ws.Range("D12") '32.04'
$ws.Range("E12").Value = '  -2.36%  '

$ws.Range("E13").Value = '  -0.58%  '

Set-TextValue $ws.Range("D14") '7.08'
$ws.Range("E14").Value = '  -6.44%  '

Set-TextValue $ws.Range("D15") '2.683.36'
$ws.Range("E15").Value = '  -1.41%  '

Set-TextValue $ws.Range("D16") '16.05'
$ws.Range("E16").Value = '  -5.01%  '

Set-TextValue $ws.Range("D17") '0.890'
$ws.Range("E17").Value = '  -3.28%  '

Set-TextValue $ws.Range("D18") '2.332.73'
$ws.Range("E18").Value = '  -1.59%  '

Set-TextValue $ws.Range("D19") '43.592.73'
$ws.Range("E19").Value = '  -0.65%  '

Set-TextValue $ws.Range("D20") '0.0000100'
$ws.Range("E20").Value = '  -2.42%  '

Set-TextValue $ws.Range("D21") '77.93'
$ws.Range("E21").Value = '  +0.29%  '

Set-TextValue $ws.Range("D22") '6.62'
$ws.Range("E22").Value = '  -0.78%  '

Set-TextValue $ws.Range("D23") '250.59'
$ws.Range("E23").Value = '  -2.67%  '

Set-TextValue $ws.Range("D25") '1.88'
$ws.Range("E25").Value = '  +7.35%  '

$ws.Range("E26").Value = '  +2.53%  '

Set-TextValue $ws.Range("D27") '2.47'
$ws.Range("E27").Value = '  -2.10%  '

Set-TextValue $ws.Range("D28") '10.26'
$ws.Range("E28").Value = '  -8.08%  '

$ws.Range("E29").Value = '  -1.03%  '

Set-TextValue $ws.Range("D30") '174.83'
$ws.Range("E30").Value = '  -0.47%  '

Set-TextValue $ws.Range("D31") '22.02'
$ws.Range("E31").Value = '  -4.72%  '

Set-TextValue $ws.Range("D32") '0.126'
$ws.Range("E32").Value = '  -2.44%  '

$ws.Range("E33").Value = '  -0.67%  '

Set-TextValue $ws.Range("D34") '0.0727'
$ws.Range("E34").Value = '  -2.92%  '

Set-TextValue $ws.Range("D35") '5.04'
$ws.Range("E35").Value = '  -4.95%  '

Set-TextValue $ws.Range("D36") '5.32'
$ws.Range("E36").Value = '  -0.36%  '

$ws.Range("E37").Value = '  -2.34%  '

$ws.Range("E38").Value = '  -3.77%  '

Set-TextValue $ws.Range("D39") '6.32'

$ws.Range("B40").Value = 'FTXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range("D40") '5.44'
$ws.Range("E40").Value = '  +20.85%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D41") '0.0268'
$ws.Range("E41").Value = '  -2.90%  '

Set-TextValue $ws.Range("D42") '64.60'
$ws.Range("E42").Value = '  +18.61%  '

$ws.Range("E43").Value = '  +2.34%  '

$ws.Range("E44").Value = '  +4.76%  '

Set-TextValue $ws.Range("D45") '18.76'
$ws.Range("E45").Value = '  -1.01%  '

Set-TextValue $ws.Range("D46") '0.194'
$ws.Range("E46").Value = '  -4.28%  '

$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D48") '2.42'
$ws.Range("E48").Value = '  -3.29%  '

$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D49") '1.21'
$ws.Range("E49").Value = '  -3.81%  '

$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D50") '2.90'
$ws.Range("E50").Value = '  +3.63%  '

$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D51") '1.14'
$ws.Range("E51").Value = '  -4.82%  '
